# The regenerated test fixture only differs from the original in the
# content of cell A2 on Sheet1: it used to hold the placeholder text "-"
# and is now blank (an empty string). Everything else in the diff
# (shared-strings -> inline-strings conversion, row "spans" attributes,
# styles.xml color-table / attribute-order churn, workbookView default
# attributes, pageMargins attribute order, ...) is just incidental
# re-serialization noise from whatever tool regenerated the fixture, not
# a deliberate spreadsheet edit, so the only actual change to make here
# is clearing A2.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = ""
